$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Text -eq "T/R1") {
        $cell.Value = "T"
    }
}
